$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per diff
$ws.Range("B4").Value = 402
$ws.Range("B6").Value = 397

# Update the active cell selection to B4
$ws.Range("B4").Select()
